# Add a new syllabus row ("named_parameters.py") under the "Functions" topic,
# and extend the merged "Functions" label to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row just below "lambda_function.py" (row 57), i.e. at row 58.
#     Everything at/after row 58 shifts down by one.
$ws.Rows("58:58").Insert()

# --- 2. Give the new row the same look (borders/font/alignment) as the row above it
#     (row 57, the last row of the old "Functions" merge block) before we grow the merge.
$ws.Range("A57:B57").Copy()
$ws.Range("A58:B58").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Put the new program name in column B.
$ws.Range("B58").Value2 = "named_parameters.py"

# --- 4. Grow the "Functions" category merge (was A56:A57) to include the new row.
$ws.Range("A56:A58").Merge()

# --- 5. Re-apply the freeze pane (top row frozen) now that the sheet has grown,
#     and leave the selection on the newly added cell.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("B58").Select()

Write-Output "Inserted named_parameters.py row and extended Functions merge"
